# Bond dates update: shift the "today" reference date forward by 2 days.
# Column G = "Dni od poprzedniej wypłaty" (days since previous payment) -> increases by 2
# Column I = "Dni do następnej wypłaty" (days until next payment)       -> decreases by 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $gCell = $ws.Cells.Item($r, 7)
    $gText = $gCell.Text
    if ($gText -ne "") {
        $gNum = [int]$gText
        $gCell.Value = $gNum + 2
    }

    $iCell = $ws.Cells.Item($r, 9)
    $iText = $iCell.Text
    if ($iText -ne "") {
        $iNum = [int]$iText
        $iCell.Value = $iNum - 2
    }
}
